$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column B value (Taxonsorteringsordning) changes from 79244 to 79245
$rows = @(2,3,4,5,6,9,10,11,12,13,15)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = 79245
}
